$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (the orphan "chốt đề tài" label row under section II),
# which shifts all subsequent rows up by one.
$ws.Rows("13").Delete()

# Update the active cell selection to match the new layout.
$ws.Range("B16").Select()
